$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric strings in columns D and E keep their exact
# formatting (trailing zeros, thousands separators, etc.) by forcing the
# cell number format to Text before assigning the value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.111.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.483.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.41%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "417.70"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.49"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.641"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.37%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.786"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +9.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +23.19%  "

$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.02"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.14%  "

$ws.Range("B12").Value = "ShibaInu"
$ws.Range("C12").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000266"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +29.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.82"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +9.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.031.54"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.47%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.42"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.471.08"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.10"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.40"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "63.988.58"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "458.44"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.86"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.27"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.48"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.21"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +12.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.33"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.21%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.78"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.42%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.81"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.93%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.59"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +7.49%  "

$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.62"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.52%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.114"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.64%  "

$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.166"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.22"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.40%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.02"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0516"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +7.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.144"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +8.27%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.11"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +7.53%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0658"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +66.09%  "

$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.74"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +9.32%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +9.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.36"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.320"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.69"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.02"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.96%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.91"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.94"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.140"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.98%  "

Write-Host "Updated cryptos list"